$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Clear / update individual cells in rows 2-25 ---
$ws.Range("D2").ClearContents()

$ws.Range("E4").ClearContents()

$ws.Range("C6").Value = 15.1
$ws.Range("E6").Value = -5.7

$ws.Range("C8").ClearContents()

$ws.Range("E12").ClearContents()

$ws.Range("E14").Value = -5.4

$ws.Range("C18").Value = 11.5

$ws.Range("C20").ClearContents()

$ws.Range("E21").Value = -8.699999999999999

$ws.Range("E22").Value = -6.1

$ws.Range("C23").Value = 12.2

$ws.Range("C25").ClearContents()

# --- Remove rows "RM 232" (row 26) and "SC 92" (originally row 28, becomes
#     row 27 once row 26 is removed) so the remaining rows shift up ---
$ws.Rows.Item(26).Delete()
$ws.Rows.Item(27).Delete()

# --- Apply the remaining per-cell edits to the now-shifted rows 26-33 ---
$ws.Range("E26").ClearContents()

$ws.Range("B27").Value = -20.4
$ws.Range("E27").ClearContents()

$ws.Range("B28").ClearContents()
$ws.Range("E28").ClearContents()

$ws.Range("B29").ClearContents()

$ws.Range("B30").Value = -19.7
$ws.Range("C30").Value = 11.4
$ws.Range("D30").Value = -13.6

$ws.Range("E31").Value = -8.1

$ws.Range("B32").ClearContents()
